$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 1303.6666
$ws.Range("I20").Value = 1303.6666
$ws.Range("K20").Value = 1303.6666
$ws.Range("M20").Value = -1073.6666
# Row 28
$ws.Range("H28").Value = 478.30768
$ws.Range("I28").Value = 484.83334
$ws.Range("K28").Value = 484.83334
$ws.Range("M28").Value = 0.1666599999999789
# Row 35
$ws.Range("H35").Value = 1303.6666
$ws.Range("I35").Value = 1303.6666
$ws.Range("K35").Value = 1303.6666
$ws.Range("M35").Value = -924.6666
# Row 132
$ws.Range("H132").Value = 3363.1667
$ws.Range("I132").Value = 3135.8
$ws.Range("K132").Value = 9407.400000000001
$ws.Range("M132").Value = -6877.400000000001
# Row 137
$ws.Range("H137").Value = 3492.3076
$ws.Range("I137").Value = 3678.625
$ws.Range("K137").Value = 11035.875
$ws.Range("M137").Value = -8485.875
# Row 138
$ws.Range("H138").Value = 2526.6956
$ws.Range("I138").Value = 1370.5333
$ws.Range("K138").Value = 4111.5999
$ws.Range("M138").Value = 1028.4001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6208.684
$ws.Range("I61").Value = 6241.4053
$ws.Range("J61").Value = 4998
$ws.Range("K61").Value = 6241.4053
$ws.Range("L61").Value = 4998
$ws.Range("M61").Value = -6029.4053
$ws.Range("N61").Value = -5422
# Row 108
$ws.Range("H108").Value = 90000
$ws.Range("J108").Value = 90000
$ws.Range("L108").Value = 90000
$ws.Range("N108").Value = -97680
# Row 110
$ws.Range("H110").Value = 4071.2856
$ws.Range("I110").Value = 2210.8
$ws.Range("K110").Value = 2210.8
$ws.Range("M110").Value = -165.8000000000002
# Row 136
$ws.Range("H136").Value = 6208.684
$ws.Range("I136").Value = 6241.4053
$ws.Range("J136").Value = 4998
$ws.Range("K136").Value = 18724.2159
$ws.Range("L136").Value = 14994
$ws.Range("M136").Value = -16174.2159
$ws.Range("N136").Value = -20094

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1781.1724
$ws.Range("I86").Value = 1523.6666
$ws.Range("J86").Value = 3017.2
$ws.Range("K86").Value = 1523.6666
$ws.Range("L86").Value = 3017.2
$ws.Range("M86").Value = -400.6666
$ws.Range("N86").Value = -5263.2
# Row 89
$ws.Range("H89").Value = 1781.1724
$ws.Range("I89").Value = 1523.6666
$ws.Range("J89").Value = 3017.2
$ws.Range("K89").Value = 7618.333000000001
$ws.Range("L89").Value = 15086
$ws.Range("M89").Value = -2002.333000000001
$ws.Range("N89").Value = -26318
# Row 134
$ws.Range("H134").Value = 5224.382
$ws.Range("I134").Value = 5032.5957
$ws.Range("K134").Value = 15097.7871
$ws.Range("M134").Value = -12562.7871

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 826
$ws.Range("I6").Value = 1175
$ws.Range("J6").Value = 477
$ws.Range("K6").Value = 1175
$ws.Range("L6").Value = 477
$ws.Range("M6").Value = -1062
$ws.Range("N6").Value = -703
# Row 109
$ws.Range("H109").Value = 30095
$ws.Range("J109").Value = 30095
$ws.Range("L109").Value = 30095
$ws.Range("N109").Value = -32175
# Row 132
$ws.Range("H132").Value = 8113.5713
$ws.Range("I132").Value = 5698.75
$ws.Range("J132").Value = 11333.333
$ws.Range("K132").Value = 17096.25
$ws.Range("L132").Value = 33999.999
$ws.Range("M132").Value = -14566.25
$ws.Range("N132").Value = -39059.999
# Row 134
$ws.Range("H134").Value = 9099.888999999999
$ws.Range("I134").Value = 6875
$ws.Range("J134").Value = 10879.8
$ws.Range("K134").Value = 20625
$ws.Range("L134").Value = 32639.4
$ws.Range("M134").Value = -18090
$ws.Range("N134").Value = -37709.39999999999
# Row 141
$ws.Range("H141").Value = 249032.78
$ws.Range("J141").Value = 249032.78
$ws.Range("L141").Value = 249032.78
$ws.Range("N141").Value = -259392.78

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 1111397.2
$ws.Range("I7").Value = 1111397.2
$ws.Range("K7").Value = 3334191.6
$ws.Range("M7").Value = -3334079.6
# Row 87
$ws.Range("H87").Value = 9874.75
$ws.Range("I87").Value = 7714
$ws.Range("K87").Value = 23142
$ws.Range("M87").Value = -21894
# Row 90
$ws.Range("H90").Value = 9874.75
$ws.Range("I90").Value = 7714
$ws.Range("K90").Value = 69426
$ws.Range("M90").Value = -63186
# Row 134
$ws.Range("H134").Value = 1149
$ws.Range("I134").Value = 1063.9
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 3191.7
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = 1878.3
$ws.Range("N134").Value = -16140

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 10000
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("N46").Value = -10312
# Row 99
$ws.Range("H99").Value = 8864.299999999999
$ws.Range("I99").Value = 5404.778
$ws.Range("J99").Value = 40000
$ws.Range("K99").Value = 5404.778
$ws.Range("L99").Value = 40000
$ws.Range("M99").Value = -3158.778
$ws.Range("N99").Value = -44492
# Row 102
$ws.Range("H102").Value = 2798.037
$ws.Range("I102").Value = 1897.25
$ws.Range("J102").Value = 10004.333
$ws.Range("K102").Value = 1897.25
$ws.Range("L102").Value = 10004.333
$ws.Range("M102").Value = -275.25
$ws.Range("N102").Value = -13248.333
# Row 132
$ws.Range("H132").Value = 4333.8335
$ws.Range("I132").Value = 3744.7856
$ws.Range("K132").Value = 11234.3568
$ws.Range("M132").Value = -8704.356800000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
# Row 40
$ws.Range("H40").Value = 4870.2144
$ws.Range("I40").Value = 4797
$ws.Range("J40").Value = 4890.1816
$ws.Range("K40").Value = 4797
$ws.Range("L40").Value = 4890.1816
$ws.Range("M40").Value = -4661
$ws.Range("N40").Value = -5162.1816
# Row 136
$ws.Range("H136").Value = 68973200
$ws.Range("I136").Value = 45462410
$ws.Range("K136").Value = 136387230
$ws.Range("M136").Value = -136384680

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4450.5
$ws.Range("I81").Value = 4450.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 8901
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -7840
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 4450.5
$ws.Range("I84").Value = 4450.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 44505
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -39201
$ws.Range("N84").ClearContents()
# Row 126
$ws.Range("H126").Value = 6159.913
$ws.Range("I126").Value = 4241.4375
$ws.Range("K126").Value = 12724.3125
$ws.Range("M126").Value = -10254.3125
# Row 132
$ws.Range("H132").Value = 5620.7188
$ws.Range("I132").Value = 4852.5
$ws.Range("K132").Value = 14557.5
$ws.Range("M132").Value = -12027.5
# Row 136
$ws.Range("H136").Value = 7753.304
$ws.Range("I136").Value = 7407.222
$ws.Range("K136").Value = 22221.666
$ws.Range("M136").Value = -19671.666
